# The document has two BTec "Logo-Orange" pictures living in the page
# headers (currently saved internally as image1.jpg) and two Pearson
# "PearsonLogo.png" pictures living in the page footers (currently saved
# internally as image2.png). The authors swapped the stored file-name
# metadata so that the BTec logos now read image2.jpg and the Pearson
# logos now read image1.png.
#
# Renaming an InlineShape is done by assigning to its .Name property.
# For header pictures this works directly off the InlineShapes
# collection; for footer pictures this runtime only applies the rename
# reliably once the picture's range has been selected first, so every
# shape is selected before its Name is updated (this also matches how
# a person would rename a picture interactively from the Selection
# pane).

$d = $word.ActiveDocument

function Rename-InlinePicture($shape, $newName) {
    $shape.Range.Select()
    $word.Selection.InlineShapes.Item(1).Name = $newName
}

$sec = $d.Sections.Item(1)

# Header pictures (BTec_Logo-Orange): image1.jpg -> image2.jpg
for ($i = 1; $i -le 3; $i++) {
    $hdr = $sec.Headers.Item($i)
    if ($hdr.Exists -and $hdr.Range.InlineShapes.Count -gt 0) {
        for ($j = 1; $j -le $hdr.Range.InlineShapes.Count; $j++) {
            Rename-InlinePicture $hdr.Range.InlineShapes.Item($j) "image2.jpg"
        }
    }
}

# Footer pictures (PearsonLogo): image2.png -> image1.png
for ($i = 1; $i -le 3; $i++) {
    $ftr = $sec.Footers.Item($i)
    if ($ftr.Exists -and $ftr.Range.InlineShapes.Count -gt 0) {
        for ($j = 1; $j -le $ftr.Range.InlineShapes.Count; $j++) {
            Rename-InlinePicture $ftr.Range.InlineShapes.Item($j) "image1.png"
        }
    }
}
